$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-27 Sunday" "2025-07-28 Monday"

Replace-Text "20×51=" "39×96="
Replace-Text "29×95=" "23×21="
Replace-Text "47×73=" "38×97="
Replace-Text "51×62=" "37×94="
Replace-Text "14×74=" "77×82="
Replace-Text "30×37=" "20×57="
Replace-Text "79×69=" "49×26="
Replace-Text "44×17=" "32×31="
Replace-Text "56×43=" "63×93="
Replace-Text "57×14=" "47×45="
Replace-Text "63×30=" "34×37="
Replace-Text "87×47=" "28×98="
Replace-Text "46×17=" "31×21="
Replace-Text "16×63=" "97×96="
Replace-Text "31×34=" "48×42="
Replace-Text "77×97=" "78×80="
Replace-Text "68×23=" "89×50="
Replace-Text "87×17=" "88×52="
Replace-Text "76×55=" "56×61="
Replace-Text "29×55=" "40×20="
Replace-Text "49×40=" "19×78="
Replace-Text "65×50=" "26×62="
Replace-Text "75×64=" "66×35="
Replace-Text "51×48=" "91×30="
Replace-Text "20×73=" "99×45="
